$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from F1 (bold/centered) onto the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New header values
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# New data cells
$ws.Range("G2").Value = 0.1260932844166139
$ws.Range("H2").Value = 0.991
